$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Adcy1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 151.7260716666667
$ws.Cells.Item(2, 8).Value = 455.178215
$ws.Cells.Item(2, 9).Value = 0.2700739458961593
$ws.Cells.Item(2, 10).Value = 0.2783366498663096
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.04069533333333333
$ws.Cells.Item(2, 14).Value = 0.122086
$ws.Cells.Item(2, 15).Value = 0.1372091977792263
$ws.Cells.Item(2, 16).Value = 0.1776788780933641
$ws.Cells.Item(2, 17).Value = 6.174543061832223
$ws.Cells.Item(2, 18).Value = 55.57088755649
$ws.Cells.Item(2, 19).Value = 0.03705662945748218
$ws.Cells.Item(2, 20).Value = 0.04945454368051139

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Adcy1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 151.7260716666667
$ws.Cells.Item(3, 8).Value = 455.178215
$ws.Cells.Item(3, 9).Value = 0.2700739458961593
$ws.Cells.Item(3, 10).Value = 0.2783366498663096
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.053234
$ws.Cells.Item(3, 14).Value = 0.159702
$ws.Cells.Item(3, 15).Value = 0.1794848164714873
$ws.Cells.Item(3, 16).Value = 0.2324236373479878
$ws.Cells.Item(3, 17).Value = 8.076985699103334
$ws.Cells.Item(3, 18).Value = 72.69287129193
$ws.Cells.Item(3, 19).Value = 0.04847417261290254
$ws.Cells.Item(3, 20).Value = 0.064692016569181

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Adcy1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 151.7260716666667
$ws.Cells.Item(4, 8).Value = 455.178215
$ws.Cells.Item(4, 9).Value = 0.2700739458961593
$ws.Cells.Item(4, 10).Value = 0.2783366498663096
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.202664
$ws.Cells.Item(4, 14).Value = 0.405328
$ws.Cells.Item(4, 15).Value = 0.6833059857492862
$ws.Cells.Item(4, 16).Value = 0.589897484558648
$ws.Cells.Item(4, 17).Value = 30.74941258825334
$ws.Cells.Item(4, 18).Value = 184.49647552952
$ws.Cells.Item(4, 19).Value = 0.1845431438257745
$ws.Cells.Item(4, 20).Value = 0.1641900896166172

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Adcy1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 82.248871
$ws.Cells.Item(5, 8).Value = 246.746613
$ws.Cells.Item(5, 9).Value = 0.146403824289839
$ws.Cells.Item(5, 10).Value = 0.150882936320401
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.04069533333333333
$ws.Cells.Item(5, 14).Value = 0.122086
$ws.Cells.Item(5, 15).Value = 0.1372091977792263
$ws.Cells.Item(5, 16).Value = 0.1776788780933641
$ws.Cells.Item(5, 17).Value = 3.347145221635333
$ws.Cells.Item(5, 18).Value = 30.124306994718
$ws.Cells.Item(5, 19).Value = 0.02008795128261961
$ws.Cells.Item(5, 20).Value = 0.02680871084884135

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Adcy1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 82.248871
$ws.Cells.Item(6, 8).Value = 246.746613
$ws.Cells.Item(6, 9).Value = 0.146403824289839
$ws.Cells.Item(6, 10).Value = 0.150882936320401
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.053234
$ws.Cells.Item(6, 14).Value = 0.159702
$ws.Cells.Item(6, 15).Value = 0.1794848164714873
$ws.Cells.Item(6, 16).Value = 0.2324236373479878
$ws.Cells.Item(6, 17).Value = 4.378436398814
$ws.Cells.Item(6, 18).Value = 39.40592758932601
$ws.Cells.Item(6, 19).Value = 0.02627726353338563
$ws.Cells.Item(6, 20).Value = 0.03506876087333242

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Adcy1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 82.248871
$ws.Cells.Item(7, 8).Value = 246.746613
$ws.Cells.Item(7, 9).Value = 0.146403824289839
$ws.Cells.Item(7, 10).Value = 0.150882936320401
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.202664
$ws.Cells.Item(7, 14).Value = 0.405328
$ws.Cells.Item(7, 15).Value = 0.6833059857492862
$ws.Cells.Item(7, 16).Value = 0.589897484558648
$ws.Cells.Item(7, 17).Value = 16.668885192344
$ws.Cells.Item(7, 18).Value = 100.013311154064
$ws.Cells.Item(7, 19).Value = 0.1000386094738337
$ws.Cells.Item(7, 20).Value = 0.08900546459822721

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Adcy1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 123.444321
$ws.Cells.Item(8, 8).Value = 370.332963
$ws.Cells.Item(8, 9).Value = 0.2197321429647646
$ws.Cells.Item(8, 10).Value = 0.2264546783208506
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.04069533333333333
$ws.Cells.Item(8, 14).Value = 0.122086
$ws.Cells.Item(8, 15).Value = 0.1372091977792263
$ws.Cells.Item(8, 16).Value = 0.1776788780933641
$ws.Cells.Item(8, 17).Value = 5.023607791202
$ws.Cells.Item(8, 18).Value = 45.212470120818
$ws.Cells.Item(8, 19).Value = 0.03014927106250561
$ws.Cells.Item(8, 20).Value = 0.0402362131830424

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Adcy1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 123.444321
$ws.Cells.Item(9, 8).Value = 370.332963
$ws.Cells.Item(9, 9).Value = 0.2197321429647646
$ws.Cells.Item(9, 10).Value = 0.2264546783208506
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.053234
$ws.Cells.Item(9, 14).Value = 0.159702
$ws.Cells.Item(9, 15).Value = 0.1794848164714873
$ws.Cells.Item(9, 16).Value = 0.2324236373479878
$ws.Cells.Item(9, 17).Value = 6.571434984114001
$ws.Cells.Item(9, 18).Value = 59.14291485702601
$ws.Cells.Item(9, 19).Value = 0.03943858335291738
$ws.Cells.Item(9, 20).Value = 0.05263342002980061

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "Adcy1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 123.444321
$ws.Cells.Item(10, 8).Value = 370.332963
$ws.Cells.Item(10, 9).Value = 0.2197321429647646
$ws.Cells.Item(10, 10).Value = 0.2264546783208506
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.202664
$ws.Cells.Item(10, 14).Value = 0.405328
$ws.Cells.Item(10, 15).Value = 0.6833059857492862
$ws.Cells.Item(10, 16).Value = 0.589897484558648
$ws.Cells.Item(10, 17).Value = 25.017719871144
$ws.Cells.Item(10, 18).Value = 150.106319226864
$ws.Cells.Item(10, 19).Value = 0.1501442885493416
$ws.Cells.Item(10, 20).Value = 0.1335850451080076

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Gnai2"
$ws.Cells.Item(11, 3).Value = "Adcy1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 154.3429766666667
$ws.Cells.Item(11, 8).Value = 463.02893
$ws.Cells.Item(11, 9).Value = 0.2747320633285943
$ws.Cells.Item(11, 10).Value = 0.2831372788071194
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.04069533333333333
$ws.Cells.Item(11, 14).Value = 0.122086
$ws.Cells.Item(11, 15).Value = 0.1372091977792263
$ws.Cells.Item(11, 16).Value = 0.1776788780933641
$ws.Cells.Item(11, 17).Value = 6.281038883108889
$ws.Cells.Item(11, 18).Value = 56.52934994798
$ws.Cells.Item(11, 19).Value = 0.03769576601354802
$ws.Cells.Item(11, 20).Value = 0.05030751404485703

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Gnai2"
$ws.Cells.Item(12, 3).Value = "Adcy1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 154.3429766666667
$ws.Cells.Item(12, 8).Value = 463.02893
$ws.Cells.Item(12, 9).Value = 0.2747320633285943
$ws.Cells.Item(12, 10).Value = 0.2831372788071194
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.053234
$ws.Cells.Item(12, 14).Value = 0.159702
$ws.Cells.Item(12, 15).Value = 0.1794848164714873
$ws.Cells.Item(12, 16).Value = 0.2324236373479878
$ws.Cells.Item(12, 17).Value = 8.216294019873333
$ws.Cells.Item(12, 18).Value = 73.94664617886001
$ws.Cells.Item(12, 19).Value = 0.04931023396536578
$ws.Cells.Item(12, 20).Value = 0.06580779620916204

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Gnai2"
$ws.Cells.Item(13, 3).Value = "Adcy1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 154.3429766666667
$ws.Cells.Item(13, 8).Value = 463.02893
$ws.Cells.Item(13, 9).Value = 0.2747320633285943
$ws.Cells.Item(13, 10).Value = 0.2831372788071194
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.202664
$ws.Cells.Item(13, 14).Value = 0.405328
$ws.Cells.Item(13, 15).Value = 0.6833059857492862
$ws.Cells.Item(13, 16).Value = 0.589897484558648
$ws.Cells.Item(13, 17).Value = 31.27976502317333
$ws.Cells.Item(13, 18).Value = 187.67859013904
$ws.Cells.Item(13, 19).Value = 0.1877260633496804
$ws.Cells.Item(13, 20).Value = 0.1670219685531003

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Gnai2"
$ws.Cells.Item(14, 3).Value = "Adcy1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 50.0323125
$ws.Cells.Item(14, 8).Value = 100.064625
$ws.Cells.Item(14, 9).Value = 0.08905802352064279
$ws.Cells.Item(14, 10).Value = 0.06118845668531954
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.04069533333333333
$ws.Cells.Item(14, 14).Value = 0.122086
$ws.Cells.Item(14, 15).Value = 0.1372091977792263
$ws.Cells.Item(14, 16).Value = 0.1776788780933641
$ws.Cells.Item(14, 17).Value = 2.036081634625
$ws.Cells.Item(14, 18).Value = 12.21648980775
$ws.Cells.Item(14, 19).Value = 0.01221957996307086
$ws.Cells.Item(14, 20).Value = 0.01087189633611198

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Gnai2"
$ws.Cells.Item(15, 3).Value = "Adcy1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 50.0323125
$ws.Cells.Item(15, 8).Value = 100.064625
$ws.Cells.Item(15, 9).Value = 0.08905802352064279
$ws.Cells.Item(15, 10).Value = 0.06118845668531954
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.053234
$ws.Cells.Item(15, 14).Value = 0.159702
$ws.Cells.Item(15, 15).Value = 0.1794848164714873
$ws.Cells.Item(15, 16).Value = 0.2324236373479878
$ws.Cells.Item(15, 17).Value = 2.663420123625
$ws.Cells.Item(15, 18).Value = 15.98052074175
$ws.Cells.Item(15, 19).Value = 0.01598456300691597
$ws.Cells.Item(15, 20).Value = 0.01422164366651177

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Gnai2"
$ws.Cells.Item(16, 3).Value = "Adcy1"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 50.0323125
$ws.Cells.Item(16, 8).Value = 100.064625
$ws.Cells.Item(16, 9).Value = 0.08905802352064279
$ws.Cells.Item(16, 10).Value = 0.06118845668531954
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.202664
$ws.Cells.Item(16, 14).Value = 0.405328
$ws.Cells.Item(16, 15).Value = 0.6833059857492862
$ws.Cells.Item(16, 16).Value = 0.589897484558648
$ws.Cells.Item(16, 17).Value = 10.1397485805
$ws.Cells.Item(16, 18).Value = 40.558994322
$ws.Cells.Item(16, 19).Value = 0.06085388055065594
$ws.Cells.Item(16, 20).Value = 0.03609491668269579

